# HAZOP.xlsx - "modifs rapport4 + change name of tache2"
#
# Fills in the third HAZOP row (item 3 "Blackout électrique") with its
# Causes / Conséquences / Mesures à prendre text, matching the style
# already used for the two rows above it, and nudges a couple of
# presentation details (header-cell vertical alignment, last-used
# selection) to match the re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- New data for row 4 (item "3") ---------------------------------------
$ws.Range("B4").Value = "Blackout électrique"
$ws.Range("C4").Value = "Réacteur surchauffe. Arrêt de la réaction et redémarrage nécessaire"
$ws.Range("D4").Value = "Générateurs de sécurité sur le site. Disques de rupture au niveau du réacteur."

# C4/D4 become wrapped + vertically centered (like the rest of the table);
# C4 already had horizontal centering, D4 already had horizontal+vertical
# centering, so only the missing pieces need to be applied.
$ws.Range("C4:D4").WrapText = $true
$ws.Range("C4").VerticalAlignment = -4108

# --- Minor formatting / UI tweaks matching the re-save -------------------
# Header cell A1 ("/") loses its vertical centering (stays horizontally
# centered only).
$ws.Range("A1").VerticalAlignment = -4107

# Last selected cell when the file was saved.
[void]$ws.Range("D5").Select()
